$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: quality_comparison
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

# C1 gets a top+bottom thin border (no left/right)
$c1a = $ws1.Range("C1")
$c1a.Style = "Normal"
$c1a.Borders.LineStyle = 1
$c1a.Borders.Item(7).LineStyle = -4142
$c1a.Borders.Item(10).LineStyle = -4142

# D1 gets a top+bottom+right thin border (no left)
$d1a = $ws1.Range("D1")
$d1a.Style = "Normal"
$d1a.Borders.LineStyle = 1
$d1a.Borders.Item(7).LineStyle = -4142

# Anonymize the "fedcore" header label
$ws1.Range("C2").Value = "approach"

# ---------------------------------------------------------------------------
# Sheet 2: computational_comparison
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

$c1b = $ws2.Range("C1")
$c1b.Style = "Normal"
$c1b.Borders.LineStyle = 1
$c1b.Borders.Item(7).LineStyle = -4142
$c1b.Borders.Item(10).LineStyle = -4142

$d1b = $ws2.Range("D1")
$d1b.Style = "Normal"
$d1b.Borders.LineStyle = 1
$d1b.Borders.Item(7).LineStyle = -4142

# Reuse the same formatting for the second metric block (F1/G1)
$c1b.Copy()
$ws2.Range("F1").PasteSpecial(-4122)

$d1b.Copy()
$ws2.Range("G1").PasteSpecial(-4122)

# Anonymize the "fedcore" header labels
$ws2.Range("C2").Value = "approach"
$ws2.Range("F2").Value = "approach"

# Remove the stray empty inline-string cell G5
$ws2.Range("G5").ClearContents()

Write-Host "edit complete"
